# evaluateRegModel.xlsx - "created dictionary and fixed pyplot issue not being
# reset between MLR and Kmeans"
#
# The regression model's feature list (the "x" row) now also includes the
# synthetic intercept/bias column "Ones" that gets added when building the
# feature dictionary, so it is prepended to the existing comma-separated
# column list stored in B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

$ws.Range("B2").Value = "Ones, S1_GRD_4TO49, S1_GRD_5TO59, S2_GRD_3TO39, S2_GRD_4TO49, S2_GRD_5TO59, S2_GRD_6TO7, S2_BEST_GRD"

# Column B is best-fit/auto-sized to its text, and the label above just grew
# by 6 characters ("Ones, ") - widen the column to keep it fitting the text,
# same as Excel's automatic best-fit column width recalculation would do.
$ws.Columns.Item(2).ColumnWidth = 103.5
